$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - add I1 "I0" and J1 "IF", matching style of existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for columns I (I0) and J (IF)
$i0Values = @(6, 6, 6, 7, 6, 5, 4, 5)
$ifValues = @(6, 6, 6, 7, 6, 5, 4, 5)

for ($r = 2; $r -le 9; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($r, 10).Value = $ifValues[$idx]
}
